$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 18:35"

# Cell updates: country re-sorts (due to updated case counts) and refreshed
# Casos totales / Nuevos casos / Casos activos / Recuperados / Muertes hoy / Muertes figures.
$updates = @(
    ,@(4, "B", 6303288)
    ,@(4, "C", 12551)
    ,@(4, "D", 3548923)
    ,@(4, "E", 2564051)
    ,@(4, "G", 350)
    ,@(4, "H", 190314)
    ,@(5, "B", 4007502)
    ,@(5, "C", 6080)
    ,@(5, "E", 673040)
    ,@(5, "G", 158)
    ,@(5, "H", 124057)
    ,@(6, "B", 3924563)
    ,@(6, "C", 75595)
    ,@(6, "D", 3030513)
    ,@(6, "E", 825521)
    ,@(6, "G", 1043)
    ,@(6, "H", 68529)
    ,@(12, "B", 488513)
    ,@(12, "C", 8959)
    ,@(12, "G", 40)
    ,@(12, "H", 29234)
    ,@(14, "B", 416501)
    ,@(14, "C", 1762)
    ,@(14, "D", 389409)
    ,@(14, "E", 15670)
    ,@(14, "G", 78)
    ,@(14, "H", 11422)
    ,@(16, "B", 340411)
    ,@(16, "C", 1735)
    ,@(16, "G", 13)
    ,@(16, "H", 41527)
    ,@(21, "B", 274943)
    ,@(21, "C", 1642)
    ,@(21, "D", 248087)
    ,@(21, "E", 20345)
    ,@(21, "G", 49)
    ,@(21, "H", 6511)
    ,@(23, "B", 248180)
    ,@(23, "C", 789)
    ,@(23, "E", 15685)
    ,@(27, "B", 130242)
    ,@(27, "C", 319)
    ,@(27, "D", 115269)
    ,@(27, "E", 5835)
    ,@(27, "G", 3)
    ,@(27, "H", 9138)
    ,@(32, "B", 116360)
    ,@(32, "C", 903)
    ,@(32, "D", 102051)
    ,@(32, "E", 7661)
    ,@(32, "G", 29)
    ,@(32, "H", 6648)
    ,@(43, "B", 76358)
    ,@(43, "C", 714)
    ,@(43, "D", 64399)
    ,@(43, "E", 9155)
    ,@(43, "G", 14)
    ,@(43, "H", 2804)
    ,@(67, "B", 36899)
    ,@(67, "C", 167)
    ,@(67, "E", 2242)
    ,@(67, "G", 3)
    ,@(67, "H", 541)
    ,@(72, "A", "Chequia")
    ,@(72, "B", 26127)
    ,@(72, "C", 354)
    ,@(72, "D", 18663)
    ,@(72, "E", 7038)
    ,@(72, "G", 1)
    ,@(72, "H", 426)
    ,@(73, "A", "Australia")
    ,@(73, "B", 26049)
    ,@(73, "C", 126)
    ,@(73, "D", 21912)
    ,@(73, "E", 3459)
    ,@(73, "G", 15)
    ,@(73, "H", 678)
    ,@(74, "A", "El Salvador")
    ,@(74, "B", 26000)
    ,@(74, "C", 96)
    ,@(74, "D", 15119)
    ,@(74, "E", 10142)
    ,@(74, "G", 8)
    ,@(74, "H", 739)
    ,@(91, "B", 11065)
    ,@(91, "C", 31)
    ,@(91, "E", 1453)
    ,@(92, "B", 10998)
    ,@(92, "C", 241)
    ,@(92, "E", 6916)
    ,@(92, "G", 5)
    ,@(92, "H", 278)
    ,@(94, "B", 9844)
    ,@(94, "C", 116)
    ,@(94, "D", 5732)
    ,@(94, "E", 3811)
    ,@(94, "G", 5)
    ,@(94, "H", 301)
    ,@(112, "E", 314)
    ,@(112, "G", 2)
    ,@(112, "H", 94)
    ,@(118, "A", "Mozambique")
    ,@(118, "B", 4207)
    ,@(118, "C", 90)
    ,@(118, "D", 2370)
    ,@(118, "E", 1811)
    ,@(118, "G", 1)
    ,@(118, "H", 26)
    ,@(119, "A", "Tunez")
    ,@(119, "B", 4196)
    ,@(119, "C", 0)
    ,@(119, "D", 1628)
    ,@(119, "E", 2487)
    ,@(119, "G", 0)
    ,@(119, "H", 81)
    ,@(120, "A", "Eslovaquia")
    ,@(120, "B", 4163)
    ,@(120, "C", 121)
    ,@(120, "D", 2617)
    ,@(120, "E", 1509)
    ,@(120, "G", 4)
    ,@(120, "H", 37)
    ,@(121, "A", "Surinam")
    ,@(121, "B", 4149)
    ,@(121, "D", 3272)
    ,@(121, "E", 805)
    ,@(121, "H", 72)
    ,@(122, "A", "Cuba")
    ,@(122, "B", 4126)
    ,@(122, "D", 3458)
    ,@(122, "E", 570)
    ,@(122, "H", 98)
    ,@(129, "B", 3101)
    ,@(129, "C", 34)
    ,@(129, "D", 1075)
    ,@(129, "E", 1927)
    ,@(129, "G", 2)
    ,@(129, "H", 99)
    ,@(139, "A", "Jordania")
    ,@(139, "B", 2233)
    ,@(139, "C", 72)
    ,@(139, "D", 1648)
    ,@(139, "E", 570)
    ,@(139, "H", 15)
    ,@(140, "A", "Aruba")
    ,@(140, "B", 2211)
    ,@(140, "D", 934)
    ,@(140, "E", 1265)
    ,@(140, "H", 12)
    ,@(141, "A", "Guinea-Bisau")
    ,@(141, "B", 2205)
    ,@(141, "D", 1127)
    ,@(141, "E", 1044)
    ,@(141, "H", 34)
    ,@(142, "A", "Benin")
    ,@(142, "B", 2194)
    ,@(142, "D", 1738)
    ,@(142, "E", 416)
    ,@(142, "H", 40)
    ,@(148, "B", 1912)
    ,@(148, "C", 116)
    ,@(148, "E", 1022)
    ,@(182, "B", 330)
    ,@(182, "C", 11)
    ,@(182, "D", 295)
    ,@(214, "A", "Montserrat")
    ,@(214, "D", 12)
    ,@(214, "H", 1)
    ,@(215, "A", "Islas Malvinas")
    ,@(215, "D", 13)
    ,@(215, "H", 0)
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $ws.Range("$col$row").Value = $val
}
